$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (Price and Volume(1h) columns) as text,
# matching the original inline-string cell representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.12%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.29%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.696"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.32%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06217"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.44%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.750"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.64%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8522"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.00%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9154"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.04%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1400"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.61%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.55%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07092"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.05%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.96%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.23%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001541"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.79%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006181"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.94%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006013"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.17%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.37%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.44%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.166"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.49%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.90%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.082"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.62%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04249"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001201"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.42%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004082"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.36%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03936"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.42%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.14%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.07%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002211"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.53%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005162"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.54%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2904"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "114.62%"
